$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style to new columns F1:G1, then set header text
$ws.Range("A1").Copy($ws.Range("F1"))
$ws.Range("A1").Copy($ws.Range("G1"))
$ws.Range("F1").Value = "TVN"
$ws.Range("G1").Value = "CTC"

# Bulk-write full data table (rows 2-51, columns A-G) per target snapshot
$data = New-Object 'object[,]' 50,7
$data[0,0] = 39
$data[0,1] = "5:45 AM"
$data[0,2] = 1169.48
$data[0,3] = 1199.48
$data[0,4] = 21.87
$data[0,5] = 0
$data[0,6] = 2
$data[1,0] = 49
$data[1,1] = "5:50 AM"
$data[1,2] = 1467.18
$data[1,3] = 1497.18
$data[1,4] = 24.93
$data[1,5] = 0
$data[1,6] = 1
$data[2,0] = 91
$data[2,1] = "6:11 AM"
$data[2,2] = 2725.88
$data[2,3] = 2755.88
$data[2,4] = 22.16
$data[2,5] = 0
$data[2,6] = 0
$data[3,0] = 100
$data[3,1] = "6:15 AM"
$data[3,2] = 2986.3175
$data[3,3] = 3016.3175
$data[3,4] = 30.9975
$data[3,5] = 0
$data[3,6] = 1.25
$data[4,0] = 120
$data[4,1] = "6:25 AM"
$data[4,2] = 3577.16
$data[4,3] = 3607.16
$data[4,4] = 49.37
$data[4,5] = 0
$data[4,6] = 0
$data[5,0] = 151
$data[5,1] = "6:40 AM"
$data[5,2] = 4506.15
$data[5,3] = 4536.15
$data[5,4] = 33.81
$data[5,5] = 0
$data[5,6] = 1
$data[6,0] = 174
$data[6,1] = "6:52 AM"
$data[6,2] = 5193.04
$data[6,3] = 5223.04
$data[6,4] = 71.97
$data[6,5] = 0
$data[6,6] = 0
$data[7,0] = 182
$data[7,1] = "6:56 AM"
$data[7,2] = 5440.89
$data[7,3] = 5470.89
$data[7,4] = 108.74
$data[7,5] = 0
$data[7,6] = 0
$data[8,0] = 230
$data[8,1] = "7:20 AM"
$data[8,2] = 6880.786667
$data[8,3] = 6910.786667
$data[8,4] = 32.193333
$data[8,5] = 0
$data[8,6] = 1.333333
$data[9,0] = 239
$data[9,1] = "7:24 AM"
$data[9,2] = 7146.62
$data[9,3] = 7176.62
$data[9,4] = 32.95
$data[9,5] = 0
$data[9,6] = 2
$data[10,0] = 257
$data[10,1] = "7:33 AM"
$data[10,2] = 7691.65
$data[10,3] = 7721.65
$data[10,4] = 35.775
$data[10,5] = 0
$data[10,6] = 0
$data[11,0] = 275
$data[11,1] = "7:43 AM"
$data[11,2] = 8248.530000000001
$data[11,3] = 8278.530000000001
$data[11,4] = 21.35
$data[11,5] = 0
$data[11,6] = 0
$data[12,0] = 331
$data[12,1] = "8:10 AM"
$data[12,2] = 9903.9
$data[12,3] = 9933.9
$data[12,4] = 66.37
$data[12,5] = 0
$data[12,6] = 0
$data[13,0] = 524
$data[13,1] = "9:47 AM"
$data[13,2] = 15702.116667
$data[13,3] = 15732.116667
$data[13,4] = 35.09
$data[13,5] = 0
$data[13,6] = 0
$data[14,0] = 550
$data[14,1] = "10:00 A"
$data[14,2] = 16484
$data[14,3] = 16514
$data[14,4] = 26.08
$data[14,5] = 0
$data[14,6] = 1
$data[15,0] = 570
$data[15,1] = "10:10 A"
$data[15,2] = 17087.03
$data[15,3] = 17117.03
$data[15,4] = 37.05
$data[15,5] = 0
$data[15,6] = 0
$data[16,0] = 993
$data[16,1] = "1:41 PM"
$data[16,2] = 29775.52
$data[16,3] = 29805.52
$data[16,4] = 41.35
$data[16,5] = 0
$data[16,6] = 0
$data[17,0] = 1016
$data[17,1] = "1:53 PM"
$data[17,2] = 30454.73
$data[17,3] = 30484.73
$data[17,4] = 24.88
$data[17,5] = 0
$data[17,6] = 2
$data[18,0] = 1027
$data[18,1] = "1:59 PM"
$data[18,2] = 30807.17
$data[18,3] = 30837.17
$data[18,4] = 20.12
$data[18,5] = 0
$data[18,6] = 0
$data[19,0] = 1036
$data[19,1] = "2:03 PM"
$data[19,2] = 31058.49
$data[19,3] = 31088.49
$data[19,4] = 60.24
$data[19,5] = 0
$data[19,6] = 1
$data[20,0] = 1078
$data[20,1] = "2:24 PM"
$data[20,2] = 32329.13
$data[20,3] = 32359.13
$data[20,4] = 18.55
$data[20,5] = 0
$data[20,6] = 1
$data[21,0] = 1086
$data[21,1] = "2:28 PM"
$data[21,2] = 32575.77
$data[21,3] = 32605.77
$data[21,4] = 21.33
$data[21,5] = 0
$data[21,6] = 2
$data[22,0] = 1103
$data[22,1] = "2:36 PM"
$data[22,2] = 33070.88
$data[22,3] = 33100.88
$data[22,4] = 25.1
$data[22,5] = 0
$data[22,6] = 1
$data[23,0] = 1141
$data[23,1] = "2:55 PM"
$data[23,2] = 34219.87
$data[23,3] = 34249.87
$data[23,4] = 21.48
$data[23,5] = 0
$data[23,6] = 2
$data[24,0] = 1185
$data[24,1] = "3:17 PM"
$data[24,2] = 35530.52
$data[24,3] = 35560.52
$data[24,4] = 115.925
$data[24,5] = 0
$data[24,6] = 3
$data[25,0] = 1196
$data[25,1] = "3:23 PM"
$data[25,2] = 35858.715
$data[25,3] = 35888.715
$data[25,4] = 56.205
$data[25,5] = 0
$data[25,6] = 2
$data[26,0] = 1204
$data[26,1] = "3:27 PM"
$data[26,2] = 36097.916667
$data[26,3] = 36127.916667
$data[26,4] = 82.06
$data[26,5] = 0
$data[26,6] = 1.333333
$data[27,0] = 1214
$data[27,1] = "3:32 PM"
$data[27,2] = 36391.61
$data[27,3] = 36421.61
$data[27,4] = 96.52
$data[27,5] = 0
$data[27,6] = 7
$data[28,0] = 1221
$data[28,1] = "3:35 PM"
$data[28,2] = 36616.4
$data[28,3] = 36646.4
$data[28,4] = 23.47
$data[28,5] = 0
$data[28,6] = 0
$data[29,0] = 1245
$data[29,1] = "3:47 PM"
$data[29,2] = 37332.09
$data[29,3] = 37362.09
$data[29,4] = 193.12
$data[29,5] = 0
$data[29,6] = 6
$data[30,0] = 1257
$data[30,1] = "3:53 PM"
$data[30,2] = 37684.8
$data[30,3] = 37714.8
$data[30,4] = 36.82
$data[30,5] = 0
$data[30,6] = 1
$data[31,0] = 1273
$data[31,1] = "4:02 PM"
$data[31,2] = 38186.04
$data[31,3] = 38216.04
$data[31,4] = 28.77
$data[31,5] = 0
$data[31,6] = 0
$data[32,0] = 1283
$data[32,1] = "4:06 PM"
$data[32,2] = 38479.19
$data[32,3] = 38509.19
$data[32,4] = 51.24
$data[32,5] = 0
$data[32,6] = 3
$data[33,0] = 1291
$data[33,1] = "4:11 PM"
$data[33,2] = 38728.38
$data[33,3] = 38758.38
$data[33,4] = 100.57
$data[33,5] = 0
$data[33,6] = 2
$data[34,0] = 1302
$data[34,1] = "4:16 PM"
$data[34,2] = 39032.49
$data[34,3] = 39062.49
$data[34,4] = 230.5
$data[34,5] = 0
$data[34,6] = 1
$data[35,0] = 1311
$data[35,1] = "4:20 PM"
$data[35,2] = 39315.5
$data[35,3] = 39345.5
$data[35,4] = 76.02
$data[35,5] = 0
$data[35,6] = 1.5
$data[36,0] = 1318
$data[36,1] = "4:24 PM"
$data[36,2] = 39526.18
$data[36,3] = 39556.18
$data[36,4] = 111.31
$data[36,5] = 0
$data[36,6] = 2.5
$data[37,0] = 1335
$data[37,1] = "4:32 PM"
$data[37,2] = 40027.72
$data[37,3] = 40057.72
$data[37,4] = 47.705
$data[37,5] = 0
$data[37,6] = 0.5
$data[38,0] = 1345
$data[38,1] = "4:38 PM"
$data[38,2] = 40347.76
$data[38,3] = 40377.76
$data[38,4] = 94.73999999999999
$data[38,5] = 0
$data[38,6] = 1
$data[39,0] = 1357
$data[39,1] = "4:43 PM"
$data[39,2] = 40690.25
$data[39,3] = 40720.25
$data[39,4] = 163.16
$data[39,5] = 0
$data[39,6] = 5
$data[40,0] = 1372
$data[40,1] = "4:51 PM"
$data[40,2] = 41141.465
$data[40,3] = 41171.465
$data[40,4] = 98.08
$data[40,5] = 0
$data[40,6] = 5
$data[41,0] = 1391
$data[41,1] = "5:01 PM"
$data[41,2] = 41726.28
$data[41,3] = 41756.28
$data[41,4] = 40.32
$data[41,5] = 0
$data[41,6] = 2
$data[42,0] = 1412
$data[42,1] = "5:11 PM"
$data[42,2] = 42334.07
$data[42,3] = 42364.07
$data[42,4] = 101.66
$data[42,5] = 0
$data[42,6] = 2
$data[43,0] = 1427
$data[43,1] = "5:19 PM"
$data[43,2] = 42809.97
$data[43,3] = 42839.97
$data[43,4] = 161.82
$data[43,5] = 0
$data[43,6] = 9
$data[44,0] = 1441
$data[44,1] = "5:26 PM"
$data[44,2] = 43229.93
$data[44,3] = 43259.93
$data[44,4] = 24.11
$data[44,5] = 0
$data[44,6] = 1
$data[45,0] = 1451
$data[45,1] = "5:30 PM"
$data[45,2] = 43516.02
$data[45,3] = 43546.02
$data[45,4] = 88.58
$data[45,5] = 0
$data[45,6] = 2
$data[46,0] = 1479
$data[46,1] = "5:44 PM"
$data[46,2] = 44340.68
$data[46,3] = 44370.68
$data[46,4] = 49.65
$data[46,5] = 0
$data[46,6] = 3
$data[47,0] = 1492
$data[47,1] = "5:51 PM"
$data[47,2] = 44738.87
$data[47,3] = 44768.87
$data[47,4] = 22.62
$data[47,5] = 0
$data[47,6] = 0
$data[48,0] = 1522
$data[48,1] = "6:06 PM"
$data[48,2] = 45637.495
$data[48,3] = 45667.495
$data[48,4] = 33.28
$data[48,5] = 0
$data[48,6] = 0.5
$data[49,0] = 1551
$data[49,1] = "6:20 PM"
$data[49,2] = 46502.28
$data[49,3] = 46532.28
$data[49,4] = 39.27
$data[49,5] = 1.01
$data[49,6] = 1

$ws.Range("A2:G51").Value = $data

